# Updated symbol list refresh: price (column D) and volume/rank label
# (column E) updates, plus two rows (42/43) that swapped their
# coin/link/price/label values.
#
# Column D holds numeric-looking values that are stored as TEXT
# (t="inlineStr") in the source workbook, not as numbers. Assigning a
# numeric-looking string straight to Range.Value lets Excel's COM layer
# coerce it to a real number (and also drops trailing zeros, e.g.
# "240.90" -> 240.9), which would change both the stored type and the
# text. To keep these as text we momentarily force the cell to Text
# number-format before assigning the string value, then restore the
# original (General) number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
}

# --- Column D: price updates (keep as text) ---
Set-TextValue "D2"  "240.90"
Set-TextValue "D3"  "22.36"
Set-TextValue "D4"  "5.531"
Set-TextValue "D5"  "0.05587"
Set-TextValue "D7"  "6.479"
Set-TextValue "D8"  "1.092"
Set-TextValue "D9"  "0.8004"
Set-TextValue "D10" "0.1421"
Set-TextValue "D11" "0.07407"
Set-TextValue "D12" "0.03251"
Set-TextValue "D13" "0.02991"
Set-TextValue "D14" "0.09243"
Set-TextValue "D15" "0.001668"
Set-TextValue "D16" "3.263"
Set-TextValue "D17" "0.04712"
Set-TextValue "D18" "0.0005748"
Set-TextValue "D19" "0.006259"
Set-TextValue "D20" "0.001052"
Set-TextValue "D21" "0.003801"
Set-TextValue "D22" "0.0001499"
Set-TextValue "D23" "0.0004776"
Set-TextValue "D24" "3.975"
Set-TextValue "D25" "2.140"
Set-TextValue "D27" "0.1312"
Set-TextValue "D40" "0.04184"
Set-TextValue "D41" "0.007007"
Set-TextValue "D42" "0.1045"
Set-TextValue "D43" "0.002969"
Set-TextValue "D44" "0.009157"
Set-TextValue "D45" "0.00005489"
Set-TextValue "D47" "0.6798"
Set-TextValue "D48" "0.03076"

# --- Column E: rank/label text updates ---
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E40").Value = "39IDEXIDEXBestin24h"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

# --- Rows 42/43: CEJI and BKEXToken swapped places ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
